# Daily attendance processing - reorders the "Recorded By" (column G)
# contributor list on each session row so that automated/system actors
# are listed before the human recorder's e-mail address.
#
# Observed transformation (exact-string swap, not a generic reorder):
#   "backup@backdoor.com, System"          -> "System, backup@backdoor.com"
#   "backup@backdoor.com, system, System"  -> "system, System, backup@backdoor.com"
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
# Any other value (single entries, "admin@admin.com, ..." combos, etc.)
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "backup@backdoor.com, system, System" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By"
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
